$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.119.10'
$ws.Range("E2").Value = '  +0.45%  '
$ws.Range("D3").Value = '1.917.46'
$ws.Range("E3").Value = '  +2.45%  '
$ws.Range("E4").Value = '  +0.13%  '
$ws.Range("E5").Value = '  -0.11%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.001'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.11%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5070'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.36%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4061'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.85%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08327'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.47%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.115'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.87%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '42.05'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.39%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '24.26'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +5.79%  '
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.917.23'
$ws.Range("E13").Value = '  +2.79%  '
$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.411'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.88%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.247'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.58%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.003'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.16%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '92.60'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.78%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001094'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.70%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06499'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.91%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '18.48'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +3.13%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.002'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.15%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.948'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.01%  '
$ws.Range("D23").Value = '30.118.25'
$ws.Range("E23").Value = '  +0.50%  '
$ws.Range("E24").Value = '  +1.99%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.199'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.16%  '
$ws.Range("D26").Value = '2.137.00'
$ws.Range("E26").Value = '  +2.63%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '21.82'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +4.07%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '162.44'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.99%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.260'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.76%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '128.89'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.08%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.133'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +5.86%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.940'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.01%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.794'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.71%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.02443'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.09%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.305'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.61%  '
$ws.Range("E37").Value = '  +1.07%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.215'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +3.37%  '
$ws.Range("E39").Value = '  +0.05%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.6454'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.34%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '8.580'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.38%  '
$ws.Range("E42").Value = '  +1.28%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.210'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.55%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '13.30'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.46%  '
$ws.Range("B45").Value = 'Decentraland'
$ws.Range("C45").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6035'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.08%  '
$ws.Range("B46").Value = 'NEARProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.173'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +8.29%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.620'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.60%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '122.43'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.08%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.207'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.26%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.135'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.26%  '
$ws.Range("B51").Value = 'Aave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '77.49'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.42%  '
